$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.5353138446807861
$ws.Range("E2").Value = 803.1408381935562
$ws.Range("F2").Value = 0.02505393693237878
$ws.Range("G2").Value = 0.02246502765404476
$ws.Range("H2").Value = 0.02083550176233098
$ws.Range("I2").Value = 0.01969277238465683
$ws.Range("J2").Value = 0.01898122342094665
$ws.Range("K2").Value = 0.0185002006178746
$ws.Range("L2").Value = 0.01773418454330555
$ws.Range("M2").Value = 0.01746541763036902
$ws.Range("N2").Value = 0.01745413728520581
$ws.Range("O2").Value = 0.01724468264969401
$ws.Range("P2").Value = 0.01695302296185638
$ws.Range("Q2").Value = 0.01663411464264217
$ws.Range("R2").Value = 0.01640426971632253
$ws.Range("S2").Value = 0.01640426971632253
$ws.Range("T2").Value = 0.01618992998665489
$ws.Range("U2").Value = 0.01607052020018998
$ws.Range("V2").Value = 0.01591614705424874
$ws.Range("W2").Value = 0.0158294161121673
$ws.Range("X2").Value = 0.01570273532001202
$ws.Range("Y2").Value = 0.01565576682638511
# Row 3
$ws.Range("C3").Value = 0.5312526226043701
$ws.Range("E3").Value = 789.5386238847095
$ws.Range("F3").Value = 0.02559950623092262
$ws.Range("G3").Value = 0.02168771787734555
$ws.Range("H3").Value = 0.02048123767485033
$ws.Range("I3").Value = 0.01901088199968131
$ws.Range("J3").Value = 0.01841867037007175
$ws.Range("K3").Value = 0.0177010743425118
$ws.Range("L3").Value = 0.01675446221336222
$ws.Range("M3").Value = 0.01637244860416344
$ws.Range("N3").Value = 0.01617700312789199
$ws.Range("O3").Value = 0.01617700312789199
$ws.Range("P3").Value = 0.01617700312789199
$ws.Range("Q3").Value = 0.01617700312789199
$ws.Range("R3").Value = 0.01617700312789199
$ws.Range("S3").Value = 0.01602279062372938
$ws.Range("T3").Value = 0.01583664673438899
$ws.Range("U3").Value = 0.01582618634731991
$ws.Range("V3").Value = 0.0155893302092848
$ws.Range("W3").Value = 0.01548234052104845
$ws.Range("X3").Value = 0.01539061644999433
$ws.Range("Y3").Value = 0.01539061644999433
# Row 4
$ws.Range("C4").Value = 0.6836338043212891
$ws.Range("E4").Value = 793.8085322804272
$ws.Range("F4").Value = 0.02535686600970501
$ws.Range("G4").Value = 0.02200955149834611
$ws.Range("H4").Value = 0.02141151137851346
$ws.Range("I4").Value = 0.0201664722816641
$ws.Range("J4").Value = 0.01978630739698559
$ws.Range("K4").Value = 0.01953727676689875
$ws.Range("L4").Value = 0.01889019843019909
$ws.Range("M4").Value = 0.01746110512896444
$ws.Range("N4").Value = 0.01746110512896444
$ws.Range("O4").Value = 0.01699841726327672
$ws.Range("P4").Value = 0.01621235779265442
$ws.Range("Q4").Value = 0.01621235779265442
$ws.Range("R4").Value = 0.01615936388898752
$ws.Range("S4").Value = 0.0160342037633635
$ws.Range("T4").Value = 0.01593657690450387
$ws.Range("U4").Value = 0.01585205863232497
$ws.Range("V4").Value = 0.01567011607739017
$ws.Range("W4").Value = 0.01567011607739017
$ws.Range("X4").Value = 0.01560659119660802
$ws.Range("Y4").Value = 0.0154738505317822
# Row 5
$ws.Range("C5").Value = 0.5312228202819824
$ws.Range("E5").Value = 801.3703015535812
$ws.Range("F5").Value = 0.02577249053104588
$ws.Range("G5").Value = 0.02248538928952186
$ws.Range("H5").Value = 0.02125607514862765
$ws.Range("I5").Value = 0.01999601627087041
$ws.Range("J5").Value = 0.01880160943514637
$ws.Range("K5").Value = 0.01807362411507225
$ws.Range("L5").Value = 0.01747030820753502
$ws.Range("M5").Value = 0.01730619267266457
$ws.Range("N5").Value = 0.01719220274534775
$ws.Range("O5").Value = 0.01685570182261403
$ws.Range("P5").Value = 0.01685570182261403
$ws.Range("Q5").Value = 0.01671587210394608
$ws.Range("R5").Value = 0.01625546335167411
$ws.Range("S5").Value = 0.01619911783577108
$ws.Range("T5").Value = 0.01603895704417158
$ws.Range("U5").Value = 0.01588370135371283
$ws.Range("V5").Value = 0.01581194510987561
$ws.Range("W5").Value = 0.01578420799246556
$ws.Range("X5").Value = 0.01562193222198817
$ws.Range("Y5").Value = 0.01562125344159027
# Row 6
$ws.Range("C6").Value = 0.5625004768371582
$ws.Range("E6").Value = 807.5334383843237
$ws.Range("F6").Value = 0.02546738654006984
$ws.Range("G6").Value = 0.02245439111919036
$ws.Range("H6").Value = 0.02109085525892033
$ws.Range("I6").Value = 0.02022568393470204
$ws.Range("J6").Value = 0.01914361242287138
$ws.Range("K6").Value = 0.01860550993764146
$ws.Range("L6").Value = 0.01841023889518831
$ws.Range("M6").Value = 0.01795917517283499
$ws.Range("N6").Value = 0.01756371529948447
$ws.Range("O6").Value = 0.01728365863864195
$ws.Range("P6").Value = 0.01660389847926487
$ws.Range("Q6").Value = 0.01635024825988327
$ws.Range("R6").Value = 0.01635024825988327
$ws.Range("S6").Value = 0.01598842139057651
$ws.Range("T6").Value = 0.01598842139057651
$ws.Range("U6").Value = 0.01598842139057651
$ws.Range("V6").Value = 0.01597224277457512
$ws.Range("W6").Value = 0.01583377735515007
$ws.Range("X6").Value = 0.01580070985858656
$ws.Range("Y6").Value = 0.01574139256109792
# Row 7
$ws.Range("C7").Value = 0.5312933921813965
$ws.Range("E7").Value = 775.2565303326628
$ws.Range("F7").Value = 0.02608602654670392
$ws.Range("G7").Value = 0.0223985403507304
$ws.Range("H7").Value = 0.02011222203093225
$ws.Range("I7").Value = 0.01904699152736205
$ws.Range("J7").Value = 0.01886329720307902
$ws.Range("K7").Value = 0.01827745157508133
$ws.Range("L7").Value = 0.01827745157508133
$ws.Range("M7").Value = 0.0179762144906939
$ws.Range("N7").Value = 0.01664710538300715
$ws.Range("O7").Value = 0.01650296802066354
$ws.Range("P7").Value = 0.01650296802066354
$ws.Range("Q7").Value = 0.01636419394695139
$ws.Range("R7").Value = 0.01595454432660061
$ws.Range("S7").Value = 0.0158678969646545
$ws.Range("T7").Value = 0.01553002158270828
$ws.Range("U7").Value = 0.01551894556156603
$ws.Range("V7").Value = 0.01537415003993115
$ws.Range("W7").Value = 0.01531354222624334
$ws.Range("X7").Value = 0.01525300206322401
$ws.Range("Y7").Value = 0.0151122130669135
# Row 8
$ws.Range("C8").Value = 0.5624563694000244
$ws.Range("E8").Value = 820.0487380709783
$ws.Range("F8").Value = 0.02563719997962122
$ws.Range("G8").Value = 0.02164900298313385
$ws.Range("H8").Value = 0.02072975242787048
$ws.Range("I8").Value = 0.02000577949120131
$ws.Range("J8").Value = 0.01930235501259281
$ws.Range("K8").Value = 0.01930235501259281
$ws.Range("L8").Value = 0.01880554607114199
$ws.Range("M8").Value = 0.01831725960957782
$ws.Range("N8").Value = 0.01773085012811446
$ws.Range("O8").Value = 0.01735293645063792
$ws.Range("P8").Value = 0.01714591294686876
$ws.Range("Q8").Value = 0.01689569087409396
$ws.Range("R8").Value = 0.01689569087409396
$ws.Range("S8").Value = 0.01673012066488883
$ws.Range("T8").Value = 0.01645162515193617
$ws.Range("U8").Value = 0.01638039478006954
$ws.Range("V8").Value = 0.01630653246013695
$ws.Range("W8").Value = 0.01622791324024736
$ws.Range("X8").Value = 0.01607319697080467
$ws.Range("Y8").Value = 0.01598535551795279
# Row 9
$ws.Range("C9").Value = 0.5312776565551758
$ws.Range("E9").Value = 784.8310951907897
$ws.Range("F9").Value = 0.02591065669243635
$ws.Range("G9").Value = 0.02223379323731333
$ws.Range("H9").Value = 0.01956148053804881
$ws.Range("I9").Value = 0.01956148053804881
$ws.Range("J9").Value = 0.01842586507020793
$ws.Range("K9").Value = 0.01825037756084482
$ws.Range("L9").Value = 0.01768605830665944
$ws.Range("M9").Value = 0.01737462389079644
$ws.Range("N9").Value = 0.01737462389079644
$ws.Range("O9").Value = 0.01686365110993961
$ws.Range("P9").Value = 0.0166245181508477
$ws.Range("Q9").Value = 0.01653245797374103
$ws.Range("R9").Value = 0.01623756744613675
$ws.Range("S9").Value = 0.01581101611768909
$ws.Range("T9").Value = 0.01581101611768909
$ws.Range("U9").Value = 0.01568520810985792
$ws.Range("V9").Value = 0.01558325028567654
$ws.Range("W9").Value = 0.01541198306736381
$ws.Range("X9").Value = 0.01541198306736381
$ws.Range("Y9").Value = 0.01529885175810506
# Row 10
$ws.Range("C10").Value = 0.5312433242797852
$ws.Range("E10").Value = 773.7485516580746
$ws.Range("F10").Value = 0.02636205806991246
$ws.Range("G10").Value = 0.02210941316145429
$ws.Range("H10").Value = 0.02093684523174621
$ws.Range("I10").Value = 0.01953744627446347
$ws.Range("J10").Value = 0.01776551551765731
$ws.Range("K10").Value = 0.01776551551765731
$ws.Range("L10").Value = 0.01713020684827979
$ws.Range("M10").Value = 0.01713020684827979
$ws.Range("N10").Value = 0.01662503900202453
$ws.Range("O10").Value = 0.01612980597521375
$ws.Range("P10").Value = 0.01612980597521375
$ws.Range("Q10").Value = 0.01568724239185356
$ws.Range("R10").Value = 0.01568724239185356
$ws.Range("S10").Value = 0.01563689501421372
$ws.Range("T10").Value = 0.01549039279099645
$ws.Range("U10").Value = 0.01530203783460348
$ws.Range("V10").Value = 0.01522993583452067
$ws.Range("W10").Value = 0.01515059117679526
$ws.Range("X10").Value = 0.01514320840494613
$ws.Range("Y10").Value = 0.01508281777111256
# Row 11
$ws.Range("C11").Value = 0.5312361717224121
$ws.Range("E11").Value = 788.5689962475135
$ws.Range("F11").Value = 0.02569595041530975
$ws.Range("G11").Value = 0.02247703006053823
$ws.Range("H11").Value = 0.02026064328427834
$ws.Range("I11").Value = 0.01912911677479032
$ws.Range("J11").Value = 0.01832675122114596
$ws.Range("K11").Value = 0.0174994554230277
$ws.Range("L11").Value = 0.01749289000651258
$ws.Range("M11").Value = 0.01688983031006103
$ws.Range("N11").Value = 0.01688983031006103
$ws.Range("O11").Value = 0.01688983031006103
$ws.Range("P11").Value = 0.01661882813474246
$ws.Range("Q11").Value = 0.0165810115908937
$ws.Range("R11").Value = 0.01633641709761944
$ws.Range("S11").Value = 0.0159077036375726
$ws.Range("T11").Value = 0.01586056416672161
$ws.Range("U11").Value = 0.0158088269644679
$ws.Range("V11").Value = 0.01563120399426205
$ws.Range("W11").Value = 0.01555776964898225
$ws.Range("X11").Value = 0.01551727940997449
$ws.Range("Y11").Value = 0.01537171532646225
